# "added names to evaluate test func to help debug"
#
# Fills in placeholder ATR (column F) values of 0 for the first 13 rows
# (the rows before the rolling ATR window in F has enough history to
# produce a "real" value), removes a stray one-off debug formula that had
# been left in G41, and updates the saved selection/view to highlight the
# newly added F1:F13 range instead of the old G42 debug cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add placeholder values (0) in column F for rows 1-13.
for ($r = 1; $r -le 13; $r++) {
    $ws.Cells.Item($r, 6).Value = 0
}

# Remove the leftover one-off debug formula in G41 (and its cached value).
$ws.Range("G41").ClearContents()

# Reflect the newly populated range in the saved selection / view instead
# of the old G42 debug-cell selection, and reset scroll position.
$ws.Range("F1:F13").Select()
